$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the PAULA row (higher row index first so lower-row deletes/inserts
# below don't shift it out from under us).
$ws.Rows.Item(510).Delete()

# Delete the WILSON row.
$ws.Rows.Item(366).Delete()

# Insert a new row for CARNEIRO right above THIAGO (row 3) and fill it in.
$ws.Rows.Item(3).Insert()
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "005685089"
$ws.Range("B3").Value = "CARNEIRO"
$ws.Range("C3").Value = 50000
